$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# multi-agent has no collision -> update arrival/collision counts and probabilities
$ws.Range("B2").Value = 40
$ws.Range("C2").Value = 0
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
